$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

# Header date
Replace-Text "2024-02-21 Wednesday" "2024-02-22 Thursday"

# Row 1
Replace-Text "323÷4=" "995÷5="
Replace-Text "747÷5=" "692÷5="
Replace-Text "678÷6=" "438÷6="
Replace-Text "571÷4=" "998÷6="
Replace-Text "706÷5=" "576÷2="

# Row 5
Replace-Text "238÷6=" "267÷2="
Replace-Text "369÷3=" "109÷2="
Replace-Text "528÷3=" "108÷6="
Replace-Text "375÷6=" "406÷2="
Replace-Text "442÷4=" "567÷4="

# Row 9 - a cell was inserted before the first cell ("696÷7=") and the
# trailing cell ("157÷9=") was removed, which (because the row keeps 5
# cells total) is equivalent to shifting every value one slot to the
# right and writing the new value into the first slot. Order matters:
# replace "557÷2=" (old cell 1) before it gets reintroduced as the new
# cell 2 value, etc.
Replace-Text "557÷2=" "696÷7="
Replace-Text "892÷5=" "557÷2="
Replace-Text "444÷4=" "922÷7="
Replace-Text "204÷6=" "928÷7="
Replace-Text "157÷9=" "842÷3="

# Row 13
Replace-Text "469÷7=" "715÷8="
Replace-Text "896÷8=" "325÷4="
Replace-Text "783÷2=" "129÷2="
Replace-Text "182÷3=" "288÷5="
Replace-Text "720÷6=" "924÷8="

# Row 17
Replace-Text "212÷9=" "483÷4="
Replace-Text "837÷2=" "511÷7="
Replace-Text "561÷3=" "415÷6="
Replace-Text "949÷7=" "447÷9="
Replace-Text "953÷9=" "301÷3="
